# Add data for 2021-12-17 (advance "through December 08" -> "through December 09")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet name and the running "through" label / header cell.
$ws.Name = "Through 2021-12-09"
$ws.Range("B1").Value2 = "December 2021 (through December 09)"

# Column B holds the current (in-progress) month's counts; the rest of the
# columns are completed prior months. New carjacking counts for 2021-12-17
# bump a handful of cells across the grid.

# Englewood
$ws.Range("B3").Value2 = 5

# North Lawndale
$ws.Range("AX4").Value2 = 2

# Garfield Park
$ws.Range("B6").Value2 = 1
$ws.Range("AL6").Value2 = 2
$ws.Range("AX6").Value2 = 3

# Austin
$ws.Range("Z7").Value2 = 2
$ws.Range("AL7").Value2 = 3

# Grand Crossing
$ws.Range("B9").Value2 = 5

# Humboldt Park
$ws.Range("AL11").Value2 = 3
$ws.Range("BJ11").Value2 = 2

# Little Italy, UIC
$ws.Range("AX12").Value2 = 1

# Roseland
$ws.Range("BJ13").Value2 = 4

# Washington Heights
$ws.Range("AX15").Value2 = 1

# United Center
$ws.Range("BJ17").Value2 = 1

# West Ridge
$ws.Range("B19").Value2 = 3

# Lake View
$ws.Range("B26").Value2 = 1

# Albany Park
$ws.Range("AX36").Value2 = 1

# Auburn Gresham
$ws.Range("B38").Value2 = 1

# Calumet Heights
$ws.Range("N40").Value2 = 3

# Pullman
$ws.Range("AL56").Value2 = 1

# Portage Park
$ws.Range("Z57").Value2 = 1

# Jefferson Park
$ws.Range("N60").Value2 = 1

# Logan Square
$ws.Range("N82").Value2 = 2
